$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.428.29"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "2.375.06"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.52%  "
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("E9").Value = "  -2.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0919"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("E14").Value = "  -3.46%  "
$ws.Range("D15").Value = "2.736.62"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("E16").Value = "  -2.90%  "
$ws.Range("D17").Value = "2.369.39"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "45.413.70"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +20.71%  "
$ws.Range("E20").Value = "  -4.33%  "
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "166.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.03%  "
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.118"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("E37").Value = "  -3.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.32%  "
$ws.Range("E39").Value = "  +7.79%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0355"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.79%  "
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.82%  "
$ws.Range("E45").Value = "  -5.79%  "
$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.25%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "1.827.54"
$ws.Range("E48").Value = "  +10.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "84.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.56%  "
